# Updated Existing MRF TC
# Update the "Date" column (B2:B61) on the NumberError sheet with new run timestamps.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NumberError")

$ws.Cells.Item(2, 2).Value = "Fri Oct 25 11:43:17 EDT 2024"
$ws.Cells.Item(3, 2).Value = "Fri Oct 25 11:43:26 EDT 2024"
$ws.Cells.Item(4, 2).Value = "Fri Oct 25 11:43:36 EDT 2024"
$ws.Cells.Item(5, 2).Value = "Fri Oct 25 11:43:45 EDT 2024"
$ws.Cells.Item(6, 2).Value = "Fri Oct 25 11:43:55 EDT 2024"
$ws.Cells.Item(7, 2).Value = "Fri Oct 25 11:44:05 EDT 2024"
$ws.Cells.Item(8, 2).Value = "Fri Oct 25 11:44:15 EDT 2024"
$ws.Cells.Item(9, 2).Value = "Fri Oct 25 11:44:24 EDT 2024"
$ws.Cells.Item(10, 2).Value = "Fri Oct 25 11:44:34 EDT 2024"
$ws.Cells.Item(11, 2).Value = "Fri Oct 25 11:44:44 EDT 2024"
$ws.Cells.Item(12, 2).Value = "Fri Oct 25 11:44:54 EDT 2024"
$ws.Cells.Item(13, 2).Value = "Fri Oct 25 11:45:04 EDT 2024"
$ws.Cells.Item(14, 2).Value = "Fri Oct 25 11:45:13 EDT 2024"
$ws.Cells.Item(15, 2).Value = "Fri Oct 25 11:45:23 EDT 2024"
$ws.Cells.Item(16, 2).Value = "Fri Oct 25 11:45:33 EDT 2024"
$ws.Cells.Item(17, 2).Value = "Fri Oct 25 11:45:43 EDT 2024"
$ws.Cells.Item(18, 2).Value = "Fri Oct 25 11:45:52 EDT 2024"
$ws.Cells.Item(19, 2).Value = "Fri Oct 25 11:46:02 EDT 2024"
$ws.Cells.Item(20, 2).Value = "Fri Oct 25 11:46:12 EDT 2024"
$ws.Cells.Item(21, 2).Value = "Fri Oct 25 11:46:21 EDT 2024"
$ws.Cells.Item(22, 2).Value = "Fri Oct 25 11:46:31 EDT 2024"
$ws.Cells.Item(23, 2).Value = "Fri Oct 25 11:46:40 EDT 2024"
$ws.Cells.Item(24, 2).Value = "Fri Oct 25 11:46:50 EDT 2024"
$ws.Cells.Item(25, 2).Value = "Fri Oct 25 11:46:59 EDT 2024"
$ws.Cells.Item(26, 2).Value = "Fri Oct 25 11:47:09 EDT 2024"
$ws.Cells.Item(27, 2).Value = "Fri Oct 25 11:47:19 EDT 2024"
$ws.Cells.Item(28, 2).Value = "Fri Oct 25 11:47:29 EDT 2024"
$ws.Cells.Item(29, 2).Value = "Fri Oct 25 11:47:38 EDT 2024"
$ws.Cells.Item(30, 2).Value = "Fri Oct 25 11:47:48 EDT 2024"
$ws.Cells.Item(31, 2).Value = "Fri Oct 25 11:47:58 EDT 2024"
$ws.Cells.Item(32, 2).Value = "Fri Oct 25 11:48:08 EDT 2024"
$ws.Cells.Item(33, 2).Value = "Fri Oct 25 11:48:18 EDT 2024"
$ws.Cells.Item(34, 2).Value = "Fri Oct 25 11:48:27 EDT 2024"
$ws.Cells.Item(35, 2).Value = "Fri Oct 25 11:48:37 EDT 2024"
$ws.Cells.Item(36, 2).Value = "Fri Oct 25 11:48:47 EDT 2024"
$ws.Cells.Item(37, 2).Value = "Fri Oct 25 11:48:57 EDT 2024"
$ws.Cells.Item(38, 2).Value = "Fri Oct 25 11:49:07 EDT 2024"
$ws.Cells.Item(39, 2).Value = "Fri Oct 25 11:49:16 EDT 2024"
$ws.Cells.Item(40, 2).Value = "Fri Oct 25 11:49:26 EDT 2024"
$ws.Cells.Item(41, 2).Value = "Fri Oct 25 11:49:35 EDT 2024"
$ws.Cells.Item(42, 2).Value = "Fri Oct 25 11:49:45 EDT 2024"
$ws.Cells.Item(43, 2).Value = "Fri Oct 25 11:49:55 EDT 2024"
$ws.Cells.Item(44, 2).Value = "Fri Oct 25 11:50:05 EDT 2024"
$ws.Cells.Item(45, 2).Value = "Fri Oct 25 11:50:15 EDT 2024"
$ws.Cells.Item(46, 2).Value = "Fri Oct 25 11:50:25 EDT 2024"
$ws.Cells.Item(47, 2).Value = "Fri Oct 25 11:50:34 EDT 2024"
$ws.Cells.Item(48, 2).Value = "Fri Oct 25 11:50:43 EDT 2024"
$ws.Cells.Item(49, 2).Value = "Fri Oct 25 11:50:53 EDT 2024"
$ws.Cells.Item(50, 2).Value = "Fri Oct 25 11:51:02 EDT 2024"
$ws.Cells.Item(51, 2).Value = "Fri Oct 25 11:51:12 EDT 2024"
$ws.Cells.Item(52, 2).Value = "Fri Oct 25 11:51:22 EDT 2024"
$ws.Cells.Item(53, 2).Value = "Fri Oct 25 11:51:32 EDT 2024"
$ws.Cells.Item(54, 2).Value = "Fri Oct 25 11:51:41 EDT 2024"
$ws.Cells.Item(55, 2).Value = "Fri Oct 25 11:51:50 EDT 2024"
$ws.Cells.Item(56, 2).Value = "Fri Oct 25 11:52:00 EDT 2024"
$ws.Cells.Item(57, 2).Value = "Fri Oct 25 11:52:09 EDT 2024"
$ws.Cells.Item(58, 2).Value = "Fri Oct 25 11:52:19 EDT 2024"
$ws.Cells.Item(59, 2).Value = "Fri Oct 25 11:52:29 EDT 2024"
$ws.Cells.Item(60, 2).Value = "Fri Oct 25 11:52:39 EDT 2024"
$ws.Cells.Item(61, 2).Value = "Fri Oct 25 11:52:48 EDT 2024"
